$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E2").Value = 65
$ws.Range("G2").Value = 78
$ws.Range("H2").Value = 96
$ws.Range("I2").Value = 105
$ws.Range("B3").Value = 73
$ws.Range("D3").Value = 118
$ws.Range("E3").Value = 125
$ws.Range("H3").Value = 131
$ws.Range("B6").Value = 340
$ws.Range("C6").Value = 428
$ws.Range("D6").Value = 369
$ws.Range("E6").Value = 400
$ws.Range("G6").Value = 409
$ws.Range("H6").Value = 409
$ws.Range("I6").Value = 455
$ws.Range("B7").Value = 460
$ws.Range("C7").Value = 571
$ws.Range("D7").Value = 577
$ws.Range("E7").Value = 600
$ws.Range("G7").Value = 616
$ws.Range("H7").Value = 650
$ws.Range("I7").Value = 762

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 5
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 5

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 8

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 4

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B7").Value = 17
$ws.Range("C8").Value = 39
$ws.Range("E16").Value = 2
$ws.Range("E19").Value = 14
$ws.Range("H19").Value = 11
$ws.Range("D20").Value = 15
$ws.Range("D21").Value = 3
$ws.Range("I32").Value = 43
$ws.Range("I36").Value = 35
$ws.Range("I41").Value = 8
$ws.Range("E53").Value = 72
$ws.Range("H53").Value = 84
$ws.Range("I53").Value = 116
$ws.Range("E65").Value = 10
$ws.Range("H65").Value = 15
$ws.Range("H74").Value = 14
$ws.Range("E76").Value = 18
$ws.Range("G76").Value = 16
$ws.Range("D77").Value = 15
$ws.Range("G80").Value = 4
$ws.Range("B88").Value = 5
$ws.Range("E94").Value = 10
$ws.Range("E95").Value = 5
$ws.Range("G96").Value = 7
$ws.Range("B98").Value = 460
$ws.Range("C98").Value = 571
$ws.Range("D98").Value = 577
$ws.Range("E98").Value = 600
$ws.Range("G98").Value = 616
$ws.Range("H98").Value = 650
$ws.Range("I98").Value = 762

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("G2").Value = 2
$ws.Range("G6").Value = 7

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H2").Value = 10
$ws.Range("H3").Value = 15
$ws.Range("E6").Value = 58
$ws.Range("I6").Value = 73
$ws.Range("E7").Value = 72
$ws.Range("H7").Value = 84
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("E6").Value = 18
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("H3").Value = 2
$ws.Range("E5").Value = 7
$ws.Range("E6").Value = 10
$ws.Range("H6").Value = 15

$ws = $wb.Worksheets.Item('River North')
$ws.Range("H5").Value = 10
$ws.Range("H6").Value = 14

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 15

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("H3").Value = 3
$ws.Range("E5").Value = 11
$ws.Range("E6").Value = 14
$ws.Range("H6").Value = 11

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("D3").Value = 2
$ws.Range("D5").Value = 2

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B3").Value = 1
$ws.Range("B6").Value = 17

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("E2").Value = 2
$ws.Range("E5").Value = 10

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("C5").Value = 30
$ws.Range("C6").Value = 39

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 5

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 3

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("D3").Value = 4
$ws.Range("D6").Value = 15
